$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: update the "Version" / date note line.
#   "Version 11.07.03, 2015-10-17"  ->  "Version 11.08.00, 2016-01-23"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Version 11.07.03, 2015-10-17", $true, $false, $false, $false, $false,
    $true, 0, $false, "Version 11.08.00, 2016-01-23", 1) | Out-Null

# ---------------------------------------------------------------------
# Change 2a: trim the sentence about surrounding the file name with
# double quotes out of the InputFile parameter description.
#   "...input file to read, surrounded by double quotes to protect
#    whitespace and special characters.    Global property..."
#   -> "...input file to read.   Global property..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    " input file to read, surrounded by double quotes to protect whitespace and special characters.    Global property values can be used with the syntax ",
    $true, $false, $false, $false, $false,
    $true, 0, $false,
    " input file to read.   Global property values can be used with the syntax ",
    1) | Out-Null

# ---------------------------------------------------------------------
# Change 2b: mention gz (gzip) support alongside zip, and restate that
# the archive holds a single compressed file.  Done in two Find/Replace
# calls so the existing "_GoBack" bookmark (which sits between the two
# runs) keeps its original position in the paragraph.
#   "...The file can be a zip file with single compressed file."
#   -> "...The file can be a zip or gz file with single compressed file."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "  The file can be a zip file with single compressed ",
    $true, $false, $false, $false, $false,
    $true, 0, $false,
    "  The file can be a zip or gz ",
    1) | Out-Null

$tail = $d.Range(900, $d.Content.End)
$tail.Find.Execute(
    "file.",
    $true, $false, $false, $false, $false,
    $true, 0, $false,
    "file with single compressed file.",
    1) | Out-Null
